# Update the dSF (column F) values for the rows that were "repulled"
# to reflect the updated mean calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -3
$ws.Range("F6").Value = -1
$ws.Range("F8").Value = -4
$ws.Range("F12").Value = -2
$ws.Range("F18").Value = -5
$ws.Range("F20").Value = -2
$ws.Range("F21").Value = -1
$ws.Range("F23").Value = 0
$ws.Range("F27").Value = 3
$ws.Range("F33").Value = -5
